# New PO forecast model
# Updates the three PO-analysis sheets (Weekly Quantity, Monthly Trend,
# PO Forecast) with the refreshed forecast data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append two new weekly rows (6 and 7)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$ws1.Range("A6").Value = 45662.99999999999
$ws1.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B6").Value = 2

$ws1.Range("A7").Value = 45669.99999999999
$ws1.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B7").Value = 7

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append one new monthly row (5)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Range("A5").Value = 45688.99999999999
$ws2.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B5").Value = 9

# ---------------------------------------------------------------------------
# Sheet 3: "PO Forecast" - refreshed forecast series: existing rows 2-13
# get updated values (and the weekly dates shift forward two weeks
# starting at row 6), plus two new forecast rows (14, 15).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

$ws3.Range("B2").Value = 88
$ws3.Range("B3").Value = 17
$ws3.Range("B4").Value = 13
$ws3.Range("B5").Value = 10

$ws3.Range("A6").Value = 45662.99999999999
$ws3.Range("B6").Value = 1

$ws3.Range("A7").Value = 45669.99999999999
$ws3.Range("A8").Value = 45676.99999999999
$ws3.Range("A9").Value = 45683.99999999999
$ws3.Range("A10").Value = 45690.99999999999
$ws3.Range("A11").Value = 45697.99999999999
$ws3.Range("A12").Value = 45704.99999999999
$ws3.Range("A13").Value = 45711.99999999999

$ws3.Range("A14").Value = 45718.99999999999
$ws3.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B14").Value = 0

$ws3.Range("A15").Value = 45725.99999999999
$ws3.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B15").Value = 0
